$wb = $excel.ActiveWorkbook

# --- "About" sheet: replace/remove notes that no longer apply now that
#     non-Kyoto-gas GWPs are hardcoded to 0 ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A10").Value = "We use GWPs for the Kyoto gases (CO2, CH4, N2O, and F-gases)"
$wsAbout.Rows("14").Delete()
$wsAbout.Rows("11:12").Delete()

# --- "GbPbT" sheet: set non-Kyoto gas GWP values (VOC, CO, NOx, PM10,
#     PM25, SOx, BC, OC) to 0 for both the 20-year and 100-year columns ---
$wsGbPbT = $wb.Worksheets.Item("GbPbT")
$wsGbPbT.Range("B3:C10").Value = 0
